$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7
$ws.Range("D2").Value = 6
$ws.Range("B4").Value = 7
$ws.Range("D4").Value = 3
